$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 51869.5
$ws.Range("I21").Value = 47850.855
$ws.Range("K21").Value = 47850.855
$ws.Range("M21").Value = -47382.855
$ws.Range("H23").Value = 51869.5
$ws.Range("I23").Value = 47850.855
$ws.Range("K23").Value = 47850.855
$ws.Range("M23").Value = -47616.855
$ws.Range("H29").Value = 6999
$ws.Range("J29").Value = 6999
$ws.Range("L29").Value = 20997
$ws.Range("N29").Value = -21559
$ws.Range("H43").Value = 5713.7646
$ws.Range("I43").Value = 4680
$ws.Range("J43").Value = 6632.6665
$ws.Range("K43").Value = 4680
$ws.Range("L43").Value = 6632.6665
$ws.Range("M43").Value = -4611
$ws.Range("N43").Value = -6770.6665
$ws.Range("H53").Value = 332
$ws.Range("I53").Value = 347.0909
$ws.Range("J53").Value = 304.33334
$ws.Range("K53").Value = 347.0909
$ws.Range("L53").Value = 304.33334
$ws.Range("M53").Value = 289.9091
$ws.Range("N53").Value = -1578.33334
$ws.Range("H112").Value = 1544.375
$ws.Range("I112").Value = 1029.6666
$ws.Range("K112").Value = 3088.9998
$ws.Range("M112").Value = -1980.9998
$ws.Range("H118").Value = 354.375
$ws.Range("I118").Value = 376.42856
$ws.Range("J118").Value = 200
$ws.Range("K118").Value = 1129.28568
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 527.71432
$ws.Range("N118").Value = -3914
$ws.Range("H128").Value = 86000
$ws.Range("J128").Value = 86000
$ws.Range("L128").Value = 86000
$ws.Range("N128").Value = -95960
$ws.Range("H135").Value = 5000
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("H137").Value = 3866.3274
$ws.Range("I137").Value = 2529.9302
$ws.Range("J137").Value = 8655.083000000001
$ws.Range("K137").Value = 7589.790599999999
$ws.Range("L137").Value = 25965.249
$ws.Range("M137").Value = -5039.790599999999
$ws.Range("N137").Value = -31065.249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16136799
$ws.Range("I32").Value = 25003548
$ws.Range("J32").Value = 15437.272
$ws.Range("K32").Value = 25003548
$ws.Range("L32").Value = 15437.272
$ws.Range("M32").Value = -25003261
$ws.Range("N32").Value = -16011.272
$ws.Range("H61").Value = 20276858
$ws.Range("I61").Value = 14712140
$ws.Range("K61").Value = 14712140
$ws.Range("M61").Value = -14711928
$ws.Range("H122").Value = 1735.6522
$ws.Range("I122").Value = 1424.9048
$ws.Range("K122").Value = 4274.7144
$ws.Range("M122").Value = -1824.7144
$ws.Range("H132").Value = 5451.9443
$ws.Range("I132").Value = 1890.6428
$ws.Range("K132").Value = 5671.928400000001
$ws.Range("M132").Value = -3141.928400000001
$ws.Range("H136").Value = 20276858
$ws.Range("I136").Value = 14712140
$ws.Range("K136").Value = 44136420
$ws.Range("M136").Value = -44133870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 93989.664
$ws.Range("J110").Value = 93989.664
$ws.Range("L110").Value = 93989.664
$ws.Range("N110").Value = -102169.664
$ws.Range("H134").Value = 373152.16
$ws.Range("I134").Value = 2749.0908
$ws.Range("J134").Value = 2002925.6
$ws.Range("K134").Value = 8247.2724
$ws.Range("L134").Value = 6008776.800000001
$ws.Range("M134").Value = -5712.2724
$ws.Range("N134").Value = -6013846.800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 900
$ws.Range("I38").Value = 900
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 900
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -523
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -689
$ws.Range("N46").ClearContents()
$ws.Range("H50").Value = 62528.43
$ws.Range("I50").Value = 3000
$ws.Range("J50").Value = 72449.836
$ws.Range("K50").Value = 3000
$ws.Range("L50").Value = 72449.836
$ws.Range("M50").Value = -2375
$ws.Range("N50").Value = -73699.836
$ws.Range("H132").Value = 3081.0625
$ws.Range("I132").Value = 2953.1333
$ws.Range("K132").Value = 8859.3999
$ws.Range("M132").Value = -6329.3999
$ws.Range("H141").Value = 378399.88
$ws.Range("J141").Value = 413324.88
$ws.Range("L141").Value = 413324.88
$ws.Range("N141").Value = -423684.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 926.3333
$ws.Range("I14").Value = 926.3333
$ws.Range("K14").Value = 2778.9999
$ws.Range("M14").Value = -2605.9999
$ws.Range("H58").Value = 7333.3335
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 30000
$ws.Range("M58").Value = -29872
$ws.Range("H98").Value = 892.5
$ws.Range("I98").Value = 892.5
$ws.Range("K98").Value = 2677.5
$ws.Range("M98").Value = -1179.5
$ws.Range("H108").Value = 1304.6666
$ws.Range("I108").Value = 457
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 1371
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 1509
$ws.Range("N108").Value = -14760
$ws.Range("H109").Value = 2888.5
$ws.Range("I109").Value = 2888.5
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 8665.5
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -7625.5
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value = 12789.2
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 14099.111
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 42297.333
$ws.Range("M110").Value = 1090
$ws.Range("N110").Value = -50477.333
$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567
$ws.Range("H114").Value = 1119.7778
$ws.Range("J114").Value = 965
$ws.Range("L114").Value = 2895
$ws.Range("N114").Value = -9403
$ws.Range("H131").Value = 8680.366
$ws.Range("I131").Value = 312.5
$ws.Range("J131").Value = 9109.486999999999
$ws.Range("K131").Value = 937.5
$ws.Range("L131").Value = 27328.461
$ws.Range("M131").Value = 4102.5
$ws.Range("N131").Value = -37408.461
$ws.Range("H132").Value = 1882.4286
$ws.Range("J132").Value = 1949.5
$ws.Range("L132").Value = 17545.5
$ws.Range("N132").Value = -22605.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 2650
$ws.Range("I46").Value = 2650
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2650
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2494
$ws.Range("N46").ClearContents()
$ws.Range("H70").Value = 15221.556
$ws.Range("I70").Value = 17284.857
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 17284.857
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = -17014.857
$ws.Range("N70").Value = -8540
$ws.Range("H73").Value = 15221.556
$ws.Range("I73").Value = 17284.857
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 17284.857
$ws.Range("L73").Value = 8000
$ws.Range("M73").Value = -16348.857
$ws.Range("N73").Value = -9872
$ws.Range("H132").Value = 22733634
$ws.Range("I132").Value = 37040510
$ws.Range("J132").Value = 10946.471
$ws.Range("K132").Value = 111121530
$ws.Range("L132").Value = 32839.413
$ws.Range("M132").Value = -111119000
$ws.Range("N132").Value = -37899.413
$ws.Range("H134").Value = 74474.25
$ws.Range("J134").Value = 74474.25
$ws.Range("L134").Value = 223422.75
$ws.Range("N134").Value = -228492.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1493.2667
$ws.Range("J22").Value = 1500.8
$ws.Range("L22").Value = 1500.8
$ws.Range("N22").Value = -2090.8
$ws.Range("H27").Value = 1493.2667
$ws.Range("J27").Value = 1500.8
$ws.Range("L27").Value = 1500.8
$ws.Range("N27").Value = -1714.8
$ws.Range("H40").Value = 4800.4375
$ws.Range("J40").Value = 5555.5557
$ws.Range("L40").Value = 5555.5557
$ws.Range("N40").Value = -5827.5557
$ws.Range("H46").Value = 3252.4736
$ws.Range("I46").Value = 2265.9167
$ws.Range("K46").Value = 2265.9167
$ws.Range("M46").Value = -2077.9167
$ws.Range("H55").Value = 71429320
$ws.Range("J55").Value = 898.25
$ws.Range("L55").Value = 898.25
$ws.Range("N55").Value = -1244.25
$ws.Range("H61").Value = 1506.6364
$ws.Range("I61").Value = 1326.75
$ws.Range("K61").Value = 1326.75
$ws.Range("M61").Value = -1124.75
$ws.Range("H105").Value = 119000
$ws.Range("J105").Value = 119000
$ws.Range("L105").Value = 119000
$ws.Range("N105").Value = -125988
$ws.Range("H113").Value = 1506.6364
$ws.Range("I113").Value = 1326.75
$ws.Range("K113").Value = 1326.75
$ws.Range("M113").Value = 843.25
$ws.Range("H129").Value = 65333
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 65333
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 65333
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -75333
$ws.Range("H132").Value = 638561.75
$ws.Range("I132").Value = 21164.777
$ws.Range("K132").Value = 63494.33099999999
$ws.Range("M132").Value = -60964.33099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H126").Value = 1286
$ws.Range("I126").Value = 1286
$ws.Range("K126").Value = 3858
$ws.Range("M126").Value = -1388
$ws.Range("H132").Value = 258994.25
$ws.Range("I132").Value = 2562.606
$ws.Range("J132").Value = 1669368.4
$ws.Range("K132").Value = 7687.818000000001
$ws.Range("L132").Value = 5008105.199999999
$ws.Range("M132").Value = -5157.818000000001
$ws.Range("N132").Value = -5013165.199999999

